$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "288.75"
Set-TextValue "E2" "-9.75%"

Set-TextValue "D3" "40.28"
Set-TextValue "E3" "-2.52%"

Set-TextValue "D4" "5.039"
Set-TextValue "E4" "-3.93%"

Set-TextValue "D5" "0.07290"
Set-TextValue "E5" "-5.90%"

Set-TextValue "B6" "FTXToken"
Set-TextValue "C6" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D6" "1.520"
Set-TextValue "E6" "-10.15%"

Set-TextValue "B7" "MXToken"
Set-TextValue "C7" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D7" "0.9167"
Set-TextValue "E7" "-2.92%"

Set-TextValue "B8" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C8" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D8" "0.1188"
Set-TextValue "E8" "-4.48%"

Set-TextValue "B9" "WazirX"
Set-TextValue "C9" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1712"
Set-TextValue "E9" "-6.79%"

Set-TextValue "B10" "MandalaExchangeToken"
Set-TextValue "C10" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.08643"
Set-TextValue "E10" "-6.27%"

Set-TextValue "B11" "BitrueCoin"
Set-TextValue "C11" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D11" "0.04174"
Set-TextValue "E11" "-3.47%"

Set-TextValue "B12" "BitMartToken"
Set-TextValue "C12" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D12" "0.1052"
Set-TextValue "E12" "0.23%"

Set-TextValue "B13" "BitForexToken"
Set-TextValue "C13" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D13" "0.001263"
Set-TextValue "E13" "-1.96%"

Set-TextValue "B14" "TigerCash"
Set-TextValue "C14" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D14" "0.005834"
Set-TextValue "E14" "-2.06%"

Set-TextValue "B15" "LEO"
Set-TextValue "C15" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D15" "3.398"
Set-TextValue "E15" "1.63%"

Set-TextValue "B16" "GateToken"
Set-TextValue "C16" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D16" "4.273"
Set-TextValue "E16" "-1.55%"

Set-TextValue "D18" "0.3260"
Set-TextValue "E18" "-2.07%"

Set-TextValue "D19" "7.795"
Set-TextValue "E19" "0.33%"

Set-TextValue "E20" "-0.25%"

Set-TextValue "D22" "0.03850"
Set-TextValue "E22" "-4.61%"

Set-TextValue "D23" "0.001267"
Set-TextValue "E23" "0.17%"

Set-TextValue "D24" "0.003813"
Set-TextValue "E24" "-7.49%"

Set-TextValue "D25" "0.0001280"
Set-TextValue "E25" "0.71%"

Set-TextValue "D26" "0.0003722"

Set-TextValue "D38" "0.02309"
Set-TextValue "E38" "-9.58%"

Set-TextValue "D39" "0.04964"
Set-TextValue "E39" "-7.13%"

Set-TextValue "D40" "0.006445"
Set-TextValue "E40" "223.59%"

Set-TextValue "D41" "0.007677"
Set-TextValue "E41" "-1.04%"

Set-TextValue "D42" "0.1269"
Set-TextValue "E42" "-3.63%"

Set-TextValue "D43" "0.007346"
Set-TextValue "E43" "-0.22%"

Set-TextValue "D44" "0.007642"
Set-TextValue "E44" "-8.48%"

Set-TextValue "D45" "0.3125"
Set-TextValue "E45" "-1.67%"

Set-TextValue "D46" "0.00006436"
Set-TextValue "E46" "-4.20%"

Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "-0.07%"

Set-TextValue "D48" "0.2414"
Set-TextValue "E48" "20.87%"

Set-TextValue "E49" "-0.04%"

Set-TextValue "D50" "0.00002100"
Set-TextValue "E50" "-0.07%"

Set-TextValue "D51" "0.0002000"
Set-TextValue "E51" "-0.07%"
